# Applies the "Added extent reports properties" change:
#  - LoginPage (sheet1): add a new row (row 4) with locator info for the
#    "invalidCredentialsMessage" element, update selection + become the
#    active/visible tab.
#  - DashboardPage (sheet2): no longer the active/visible tab.

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("LoginPage")
$dashboardSheet = $wb.Worksheets.Item("DashboardPage")

# New data row on LoginPage. The xpath string is written first so it is
# appended to the shared-strings table before the element name (matching
# the target shared-strings ordering).
$loginSheet.Range("C4").Value = "//p[@class='oxd-text oxd-text--p oxd-alert-content-text']"
$loginSheet.Range("B4").Value = "XPATH"
$loginSheet.Range("A4").Value = "invalidCredentialsMessage"

# Update the selection on LoginPage to A6
$loginSheet.Range("A6").Select()

# Make LoginPage the active/visible sheet (mirrors tabSelected moving from
# DashboardPage to LoginPage, and activeTab clearing on the workbook view)
$loginSheet.Activate()

$wb.Save()
